# Update betting-odds values on Sheet1 to match the latest FlashScore export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (Santos vs Novorizontino)
$ws.Range("H4").Value = 3.6
$ws.Range("L4").Value = 6.5
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("W4").Value = 5
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 23
$ws.Range("AN4").Value = 3.4
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.38
$ws.Range("AU4").Value = 10
$ws.Range("AY4").Value = 41
$ws.Range("BA4").Value = 201

# Row 5 (Brusque vs Amazonas)
$ws.Range("G5").Value = 2.88
$ws.Range("H5").Value = 2.63
$ws.Range("I5").Value = 2.9
$ws.Range("J5").Value = 3.75
$ws.Range("O5").Value = 1.67
$ws.Range("P5").Value = 2.1
$ws.Range("W5").Value = 6
$ws.Range("Z5").Value = 29
$ws.Range("AH5").Value = 6.5
$ws.Range("AJ5").Value = 13
$ws.Range("AL5").Value = 34
$ws.Range("AN5").Value = 4.5
$ws.Range("AX5").Value = 19
$ws.Range("AY5").Value = 41
$ws.Range("BA5").Value = 126

# Row 8
$ws.Range("G8").Value = 3.4
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 2.3
$ws.Range("K8").Value = 1.95
$ws.Range("L8").Value = 3.1
$ws.Range("O8").Value = 1.44
$ws.Range("P8").Value = 2.63
$ws.Range("S8").Value = 1.53
$ws.Range("T8").Value = 2.38
$ws.Range("W8").Value = 8
$ws.Range("Y8").Value = 13
$ws.Range("AC8").Value = 7
$ws.Range("AK8").Value = 21
$ws.Range("AL8").Value = 21
$ws.Range("AT8").Value = 2.38
$ws.Range("AX8").Value = 13

# Row 19
$ws.Range("M19").Value = 1.07
$ws.Range("N19").Value = 9
$ws.Range("Q19").Value = 2.25
$ws.Range("R19").Value = 1.62

# Row 23
$ws.Range("M23").Value = 1.04
$ws.Range("N23").Value = 13
$ws.Range("Q23").Value = 1.85
$ws.Range("R23").Value = 2
